# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210   (columns A:J)
#   *_new -> *_FV2304   (columns L:U)
# Then wrap the data range in a real Excel Table (ListObject) with an
# AutoFilter, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" / "_new" header suffixes -----------------------
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2210')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2304')
}

# --- 2. Turn A1:U75 into a real table with an AutoFilter ------------------
$rng = $ws.Range("A1:U75")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1) -------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
